$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.117.77'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '1.813.89'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'224.43"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = "'0.555"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').Value = "'0.998"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'31.94"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.06%  '
$ws.Range('D9').Value = "'0.290"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.52%  '
$ws.Range('D10').Value = "'0.0746"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.51%  '
$ws.Range('D11').Value = "'0.0928"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '2.073.71'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').Value = '1.808.25'
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').Value = "'10.93"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').Value = "'0.642"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').Value = '34.134.14'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = "'4.32"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('D18').Value = "'69.51"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').Value = "'249.88"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').Value = '0.0₃0807'
$ws.Range('E20').Value = '  +9.03%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = "'11.05"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.30%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'0.998"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = "'4.27"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').Value = "'160.57"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.06%  '
$ws.Range('D26').Value = "'16.67"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('D27').Value = "'7.23"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.12%  '
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('D29').Value = "'0.998"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Value = "'0.0531"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.34%  '
$ws.Range('D31').Value = "'3.77"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('E32').Value = '  +1.79%  '
$ws.Range('D33').Value = "'3.57"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('D34').Value = "'1.88"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('D35').Value = '1.429.77'
$ws.Range('E35').Value = '  -1.14%  '
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('D37').Value = "'0.643"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('D38').Value = "'0.0190"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').Value = "'0.962"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.03%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = "'81.15"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = "'2.76"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.85%  '
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = "'2.16"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.60%  '
$ws.Range('D44').Value = "'6.04"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.99%  '
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').Value = "'0.0497"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.86%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.969.81'
$ws.Range('E46').Value = '  +1.29%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = "'106.54"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.51%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = "'1.04"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('D49').Value = "'0.996"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').Value = "'11.90"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.10%  '
$ws.Range('D51').Value = '0.0₆0123'
$ws.Range('E51').Value = '  +5.22%  '
